$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data for Cycle Sort (rows 12-16, columns P:U)
# Row 12 (Trial 3)
$ws.Range("P12").Value = 2
$ws.Range("Q12").Value = 9
$ws.Range("R12").Value = 277
$ws.Range("S12").Value = 26488
$ws.Range("T12").Value = "Unmeasureable"
$ws.Range("U12").Value = "Unmeasureable"

# Row 13 (Trial 4)
$ws.Range("P13").Value = 1
$ws.Range("Q13").Value = 4
$ws.Range("R13").Value = 258
$ws.Range("S13").Value = 26815
$ws.Range("T13").Value = "Unmeasureable"
$ws.Range("U13").Value = "Unmeasureable"

# Row 14 (Trial 5)
$ws.Range("P14").Value = 2
$ws.Range("Q14").Value = 3
$ws.Range("R14").Value = 233
$ws.Range("S14").Value = 23822
$ws.Range("T14").Value = "Unmeasureable"
$ws.Range("U14").Value = "Unmeasureable"

# Row 15
$ws.Range("P15").Value = 1
$ws.Range("Q15").Value = 3
$ws.Range("R15").Value = 234
$ws.Range("S15").Value = 23903
$ws.Range("T15").Value = "Unmeasureable"
$ws.Range("U15").Value = "Unmeasureable"

# Row 16 (Average row - cells already exist with style, just set values)
$ws.Range("P16").Value = 1
$ws.Range("Q16").Value = 3
$ws.Range("R16").Value = 236
$ws.Range("S16").Value = 23665
$ws.Range("T16").Value = "Unmeasureable"
$ws.Range("U16").Value = "Unmeasureable"

# Update sheet view: scroll position and selection
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("W13").Select()
